$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the adjacent header cell (G1) onto the new header cell (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the header text for the new "Save" column
$ws.Range("H1").Value = "Save"

# Values for the new "Save" column, rows 2-8
$saveValues = @(1, 0, 1, 0, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
